$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Agrp / ECs -> Mc3r / MuSCs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2239496666666667
$ws.Range("H2").Value = 0.671849
$ws.Range("I2").Value = 0.4264743968982249
$ws.Range("J2").Value = 0.4264743968982249
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01499333333333333
$ws.Range("N2").Value = 0.04498
$ws.Range("Q2").Value = 0.003357752002222222
$ws.Range("R2").Value = 0.03021976802
$ws.Range("S2").Value = 0.4264743968982249
$ws.Range("T2").Value = 0.4264743968982249

# Row 3 (Agrp / FAPs -> Mc3r / MuSCs)
$ws.Range("I3").Value = 0.4001470143891285
$ws.Range("J3").Value = 0.4001470143891285
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01499333333333333
$ws.Range("N3").Value = 0.04498
$ws.Range("Q3").Value = 0.003150469168888888
$ws.Range("R3").Value = 0.02835422252
$ws.Range("S3").Value = 0.4001470143891285
$ws.Range("T3").Value = 0.4001470143891285

# Row 4 (Agrp / MuSCs -> Mc3r / MuSCs)
$ws.Range("G4").Value = 0.09104433333333334
$ws.Range("H4").Value = 0.273133
$ws.Range("I4").Value = 0.1733785887126465
$ws.Range("J4").Value = 0.1733785887126465
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01499333333333333
$ws.Range("N4").Value = 0.04498
$ws.Range("Q4").Value = 0.001365058037777778
$ws.Range("R4").Value = 0.01228552234
$ws.Range("S4").Value = 0.1733785887126465
$ws.Range("T4").Value = 0.1733785887126465
